$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = " Quellijnstraat 56 3, 1072 XT Amsterdam Verkocht Width"
$ws.Range("C2").Value = 649000
$ws.Range("D2").Value = 67
$ws.Range("E2").Value = "A"
$ws.Range("G2").Value = 2
$ws.Range("H2").Value = 1906
$ws.Range("J2").Value = "Uitstekend"
$ws.Range("L2").Value = 0.6920833333333334
$ws.Range("B3").Value = " Eerste Sweelinckstraat 5 3, 1073 CK Amsterdam Verkocht Width"
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 70
$ws.Range("E3").Value = "A"
$ws.Range("G3").Value = 0
$ws.Range("L3").Value = 0.6833333333333333
$ws.Range("B4").Value = " Quellijnstraat 19 A, 1072 XM Amsterdam Verkocht Width"
$ws.Range("C4").Value = 725000
$ws.Range("D4").Value = 77
$ws.Range("E4").Value = "C"
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 1879
$ws.Range("J4").Value = "Goed"
$ws.Range("L4").Value = 0.6599621212121213
$ws.Range("B5").Value = " Quellijnstraat 17 3, 1072 XM Amsterdam Verkocht Width"
$ws.Range("C5").Value = 675000
$ws.Range("D5").Value = 76
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 1879
$ws.Range("L5").Value = 0.6566287878787879
$ws.Range("B6").Value = " Quellijnstraat 41 B, 1072 XP Amsterdam Verkocht Width"
$ws.Range("C6").Value = 785000
$ws.Range("D6").Value = 74
$ws.Range("E6").Value = "A"
$ws.Range("H6").Value = 1906
$ws.Range("L6").Value = 0.6554166666666668
$ws.Range("B7").Value = " Quellijnstraat 37 B, 1072 XP Amsterdam Verkocht Width"
$ws.Range("C7").Value = 725000
$ws.Range("D7").Value = 74
$ws.Range("E7").Value = "B"
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 1879
$ws.Range("L7").Value = 0.652689393939394
$ws.Range("B8").Value = " Eerste Jan van der Heijdenstraat 32 3, 1072 TV Amsterdam Verkocht Width"
$ws.Range("C8").Value = 750000
$ws.Range("D8").Value = 78
$ws.Range("E8").Value = "A"
$ws.Range("H8").Value = 1906
$ws.Range("L8").Value = 0.640625
$ws.Range("B9").Value = " Eerste Jan van der Heijdenstraat 42 1, 1072 TV Amsterdam Verkocht Width"
$ws.Range("C9").Value = 675000
$ws.Range("D9").Value = 77
$ws.Range("E9").Value = "A"
$ws.Range("L9").Value = 0.6372916666666668
$ws.Range("B10").Value = " Eerste Jan van der Heijdenstraat 42 5, 1072 TV Amsterdam Verkocht Width"
$ws.Range("C10").Value = 725000
$ws.Range("D10").Value = 79
$ws.Range("E10").Value = "A"
$ws.Range("H10").Value = 1906
$ws.Range("L10").Value = 0.6372916666666668
$ws.Range("B11").Value = " Eerste Jan van der Heijdenstraat 48 1, 1072 TW Amsterdam Verkocht Width"
$ws.Range("C11").Value = 695000
$ws.Range("D11").Value = 77
$ws.Range("L11").Value = 0.6372916666666668
$ws.Range("B12").Value = " Eerste Jan van der Heijdenstraat 44 1, 1072 TW Amsterdam Verkocht Width"
$ws.Range("C12").Value = 675000
$ws.Range("D12").Value = 77
$ws.Range("E12").Value = "A"
$ws.Range("H12").Value = 1906
$ws.Range("L12").Value = 0.6372916666666668
$ws.Range("B13").Value = " Eerste Jan van der Heijdenstraat 36 C, 1072 TV Amsterdam Verkocht Width"
$ws.Range("C13").Value = 699000
$ws.Range("D13").Value = 76
$ws.Range("E13").Value = "A"
$ws.Range("H13").Value = 1896
$ws.Range("L13").Value = 0.6339583333333334
$ws.Range("B14").Value = " Quellijnstraat 82 H, 1072 XX Amsterdam Verkocht Width"
$ws.Range("C14").Value = 425000
$ws.Range("D14").Value = 61
$ws.Range("E14").Value = "C"
$ws.Range("H14").Value = 1879
$ws.Range("L14").Value = 0.6332954545454546
$ws.Range("B15").Value = " Quellijnstraat 104 1, 1072 XZ Amsterdam Verkocht Width"
$ws.Range("C15").Value = 675000
$ws.Range("D15").Value = 80
$ws.Range("E15").Value = "C"
$ws.Range("H15").Value = 1879
$ws.Range("L15").Value = 0.6232954545454547
$ws.Range("B16").Value = " Eerste Jan van der Heijdenstraat 36 F, 1072 TV Amsterdam Verkocht Width"
$ws.Range("C16").Value = 875000
$ws.Range("D16").Value = 83
$ws.Range("E16").Value = "B"
$ws.Range("H16").Value = 1896
$ws.Range("L16").Value = 0.6212310606060606
